$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns L:S (12-19)
$headers = @(
    "hzj-混合调节_20170516_152754_ASIC_EEG",
    "hzj-混合调节_20170518_134207_ASIC_EEG",
    "hzj-混合调节_20170519_135415_ASIC_EEG",
    "zyx-混合调节_20170516_111228_ASIC_EEG",
    "zyx-混合调节_20170517_110944_ASIC_EEG",
    "zyx-混合调节_20170518_112337_ASIC_EEG",
    "zyx-混合调节_20170519_124954_ASIC_EEG",
    "zyx-混合调节_20170522_111557_ASIC_EEG"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 12 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Row 2 values for columns L:S
$row2 = @(1, 0.9862542955326461, 0.9662921348314607, 0.97435897435897434, 0.95833333333333337, 0.95145631067961167, 0.99354838709677418, 0.95145631067961167)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $col = 12 + $i
    $ws.Cells.Item(2, $col).Value = $row2[$i]
}

# Row 3 values for columns L:S
$row3 = @(1, 0.99032258064516132, 1, 0.99285714285714288, 0.97765363128491622, 0.98969072164948457, 0.96308724832214765, 0.95670995670995673)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $col = 12 + $i
    $ws.Cells.Item(3, $col).Value = $row3[$i]
}

$ws.Range("A1:S3").Select()
